$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.750.22'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').Value = '3.481.19'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '580.12'
$ws.Range('E5').Value = '  -3.23%  '
$ws.Range('D6').Value = '130.35'
$ws.Range('E6').Value = '  -3.93%  '
$ws.Range('D7').Value = '3.480.46'
$ws.Range('E7').Value = '  -1.78%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('D11').Value = '7.16'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').Value = '4.070.06'
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').Value = '27.24'
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('E16').Value = '  -3.54%  '
$ws.Range('D17').Value = '3.476.72'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').Value = '63.841.54'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D19').Value = '10.07'
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').Value = '14.30'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('E21').Value = '  -2.36%  '
$ws.Range('D22').Value = '382.99'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('D23').Value = '0.573'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').Value = '3.619.94'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').Value = '72.85'
$ws.Range('E25').Value = '  -2.48%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -4.12%  '
$ws.Range('D28').Value = '1.56'
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('E30').Value = '  -6.66%  '
$ws.Range('E31').Value = '  -3.52%  '
$ws.Range('D32').Value = '8.15'
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('D33').Value = '3.486.92'
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('D37').Value = '5.23'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').Value = '167.60'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').Value = '6.84'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').Value = '0.0795'
$ws.Range('E41').Value = '  -4.12%  '
$ws.Range('D42').Value = '26.65'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').Value = '0.810'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('D46').Value = '41.22'
$ws.Range('E46').Value = '  -3.87%  '
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('D49').Value = '2.432.22'
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').Value = '0.884'
$ws.Range('E51').Value = '  -1.54%  '
